# Update the "Yes - Only Yes" column (column 6) of the response-patterns
# table with corrected counts/percentages (includes a trailing "%" that
# was missing before, plus several recomputed values).

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$updates = @{
    3  = "37 (45.7%)"
    4  = "3 (50.0%)"
    5  = "8 (19.0%)"
    6  = "5 (11.9%)"
    7  = "3 (33.3%)"
    8  = "2 (9.1%)"
    9  = "2 (18.2%)"
    10 = "4 (21.1%)"
    11 = "0 (0.0%)"
    12 = "1 (10.0%)"
    13 = "0 (0.0%)"
    14 = "1 (25.0%)"
    15 = "0 (0.0%)"
    16 = "0 (0.0%)"
    17 = "0 (0.0%)"
    18 = "2 (100.0%)"
    19 = "9 (36.0%)"
    20 = "2 (25.0%)"
    21 = "0 (0.0%)"
    22 = "5 (38.5%)"
    23 = "3 (25.0%)"
    24 = "2 (5.6%)"
    25 = "1 (3.7%)"
    26 = "1 (5.0%)"
    27 = "0 (0.0%)"
    28 = "23 (37.7%)"
}

foreach ($row in $updates.Keys) {
    $t.Cell($row, 6).Range.Text = $updates[$row]
}
